$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1112.6897
$ws.Range("J19").Value = 1336.6842
$ws.Range("L19").Value = 1336.6842
$ws.Range("N19").Value = -1686.6842

$ws.Range("H32").Value = 1251.5
$ws.Range("J32").Value = 1251.5
$ws.Range("L32").Value = 1251.5
$ws.Range("N32").Value = -1903.5

$ws.Range("H33").Value = 785.41174
$ws.Range("I33").Value = 317.92856
$ws.Range("J33").Value = 2967
$ws.Range("K33").Value = 317.92856
$ws.Range("L33").Value = 2967
$ws.Range("M33").Value = -88.92856
$ws.Range("N33").Value = -3425

$ws.Range("H86").Value = 3184.7812
$ws.Range("I86").Value = 1408.625
$ws.Range("J86").Value = 4960.9375
$ws.Range("K86").Value = 1408.625
$ws.Range("L86").Value = 4960.9375
$ws.Range("M86").Value = -285.625
$ws.Range("N86").Value = -7206.9375

$ws.Range("H89").Value = 3184.7812
$ws.Range("I89").Value = 1408.625
$ws.Range("J89").Value = 4960.9375
$ws.Range("K89").Value = 7043.125
$ws.Range("L89").Value = 24804.6875
$ws.Range("M89").Value = -1427.125
$ws.Range("N89").Value = -36036.6875

$ws.Range("H100").Value = 1719
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 1625.3334
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 1625.3334
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -2707.3334

$ws.Range("H113").Value = 41714
$ws.Range("J113").Value = 1736
$ws.Range("L113").Value = 1736
$ws.Range("N113").Value = -8244

$ws.Range("H129").Value = 4454.143
$ws.Range("J129").Value = 1157
$ws.Range("L129").Value = 3471
$ws.Range("N129").Value = -13471

$ws.Range("H132").Value = 8071364.5
$ws.Range("I132").Value = 9266935
$ws.Range("K132").Value = 27800805
$ws.Range("M132").Value = -27798275

$ws.Range("H137").Value = 2288.2354
$ws.Range("I137").Value = 1600
$ws.Range("J137").Value = 3550
$ws.Range("K137").Value = 4800
$ws.Range("L137").Value = 10650
$ws.Range("M137").Value = -2250
$ws.Range("N137").Value = -15750

$ws.Range("H139").Value = 49977.5
$ws.Range("J139").Value = 49970
$ws.Range("L139").Value = 49970
$ws.Range("N139").Value = -60250

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 114042
$ws.Range("I2").Value = 3479.6667
$ws.Range("J2").Value = 335166.66
$ws.Range("K2").Value = 3479.6667
$ws.Range("L2").Value = 335166.66
$ws.Range("M2").Value = -3366.6667
$ws.Range("N2").Value = -335392.66

$ws.Range("H61").Value = 1599.4375
$ws.Range("I61").Value = 1381
$ws.Range("J61").Value = 2080
$ws.Range("K61").Value = 1381
$ws.Range("L61").Value = 2080
$ws.Range("M61").Value = -1169
$ws.Range("N61").Value = -2504

$ws.Range("H116").Value = 114042
$ws.Range("I116").Value = 3479.6667
$ws.Range("J116").Value = 335166.66
$ws.Range("K116").Value = 3479.6667
$ws.Range("L116").Value = 335166.66
$ws.Range("M116").Value = -1185.6667
$ws.Range("N116").Value = -339754.66

$ws.Range("H122").Value = 8000
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900

$ws.Range("H136").Value = 1599.4375
$ws.Range("I136").Value = 1381
$ws.Range("J136").Value = 2080
$ws.Range("K136").Value = 4143
$ws.Range("L136").Value = 6240
$ws.Range("M136").Value = -1593
$ws.Range("N136").Value = -11340

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 114042
$ws.Range("I3").Value = 3479.6667
$ws.Range("J3").Value = 335166.66
$ws.Range("K3").Value = 3479.6667
$ws.Range("L3").Value = 335166.66
$ws.Range("M3").Value = -3365.6667
$ws.Range("N3").Value = -335394.66

$ws.Range("H105").Value = 287738.72
$ws.Range("I105").Value = 201836
$ws.Range("J105").Value = 502495.5
$ws.Range("K105").Value = 201836
$ws.Range("L105").Value = 502495.5
$ws.Range("M105").Value = -200089
$ws.Range("N105").Value = -505989.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23704.648
$ws.Range("I31").Value = 1005.4737
$ws.Range("J31").Value = 49843.09
$ws.Range("K31").Value = 1005.4737
$ws.Range("L31").Value = 49843.09
$ws.Range("M31").Value = -710.4737
$ws.Range("N31").Value = -50433.09

$ws.Range("H34").Value = 23704.648
$ws.Range("I34").Value = 1005.4737
$ws.Range("J34").Value = 49843.09
$ws.Range("K34").Value = 1005.4737
$ws.Range("L34").Value = 49843.09
$ws.Range("M34").Value = -803.4737
$ws.Range("N34").Value = -50247.09

$ws.Range("H58").Value = 6938.6895
$ws.Range("I58").Value = 1866
$ws.Range("J58").Value = 14125
$ws.Range("K58").Value = 1866
$ws.Range("L58").Value = 14125
$ws.Range("M58").Value = -1663
$ws.Range("N58").Value = -14531

$ws.Range("H132").Value = 3982.1738
$ws.Range("I132").Value = 3598.7693
$ws.Range("J132").Value = 4480.6
$ws.Range("K132").Value = 10796.3079
$ws.Range("L132").Value = 13441.8
$ws.Range("M132").Value = -8266.3079
$ws.Range("N132").Value = -18501.8

$ws.Range("H134").Value = 1175.9722
$ws.Range("I134").Value = 1081.7742
$ws.Range("J134").Value = 1760
$ws.Range("K134").Value = 3245.3226
$ws.Range("L134").Value = 5280
$ws.Range("M134").Value = -710.3226000000004
$ws.Range("N134").Value = -10350

$ws.Range("H136").Value = 6938.6895
$ws.Range("I136").Value = 1866
$ws.Range("J136").Value = 14125
$ws.Range("K136").Value = 5598
$ws.Range("L136").Value = 42375
$ws.Range("M136").Value = -3048
$ws.Range("N136").Value = -47475

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3257
$ws.Range("I51").Value = 4900
$ws.Range("J51").Value = 2983.1667
$ws.Range("K51").Value = 14700
$ws.Range("L51").Value = 8949.500100000001
$ws.Range("M51").Value = -14240
$ws.Range("N51").Value = -9869.500100000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 11104
$ws.Range("J52").Value = 11104
$ws.Range("L52").Value = 11104
$ws.Range("N52").Value = -11622

$ws.Range("H102").Value = 354797.47
$ws.Range("I102").Value = 1811.4445
$ws.Range("J102").Value = 751906.75
$ws.Range("K102").Value = 1811.4445
$ws.Range("L102").Value = 751906.75
$ws.Range("M102").Value = -189.4445000000001
$ws.Range("N102").Value = -755150.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2549.6
$ws.Range("I68").Value = 1370.4
$ws.Range("J68").Value = 3728.8
$ws.Range("K68").Value = 1370.4
$ws.Range("L68").Value = 3728.8
$ws.Range("M68").Value = -621.4000000000001
$ws.Range("N68").Value = -5226.8

$ws.Range("H71").Value = 2549.6
$ws.Range("I71").Value = 1370.4
$ws.Range("J71").Value = 3728.8
$ws.Range("K71").Value = 6852
$ws.Range("L71").Value = 18644
$ws.Range("M71").Value = -3108
$ws.Range("N71").Value = -26132

$ws.Range("H136").Value = 1735.871
$ws.Range("I136").Value = 1645.0435
$ws.Range("J136").Value = 1997
$ws.Range("K136").Value = 4935.1305
$ws.Range("L136").Value = 5991
$ws.Range("M136").Value = -2385.1305
$ws.Range("N136").Value = -11091

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H92").Value = 19666.666
$ws.Range("J92").Value = 19666.666
$ws.Range("L92").Value = 19666.666
$ws.Range("N92").Value = -24658.666

$ws.Range("H132").Value = 2691.359
$ws.Range("I132").Value = 2538.8572
$ws.Range("J132").Value = 3079.5454
$ws.Range("K132").Value = 7616.571599999999
$ws.Range("L132").Value = 9238.636200000001
$ws.Range("M132").Value = -5086.571599999999
$ws.Range("N132").Value = -14298.6362

$ws.Range("H136").Value = 942.3570999999999
$ws.Range("I136").Value = 582.55554
$ws.Range("J136").Value = 1590
$ws.Range("K136").Value = 1747.66662
$ws.Range("L136").Value = 4770
$ws.Range("M136").Value = 802.33338
$ws.Range("N136").Value = -9870
